$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "correct "
$ws.Range("C1").Value = "not correct"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 6

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 5

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3

$ws.Range("C7").Select()
